# Updates cryptos list values (Price / Volume(1h) columns) to match the
# latest scrape. Values are stored as literal text in the sheet (mirrors
# the original inline-string cells), so cells whose new value could be
# misread as a number are forced to Text format before the write and
# restored to the default "Normal" style afterwards (no lasting format
# change, matching the unstyled source cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.174.88'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '2.448.45'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  -0.07%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '571.39'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.85%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '146.67'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("E7").Value = '  +0.10%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.537'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").Value = '2.445.08'
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("E11").Value = '  +1.24%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '5.30'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("E13").Value = '  +0.00%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '27.00'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("E15").Value = '  -1.23%  '
$ws.Range("D16").Value = '2.892.88'
$ws.Range("D17").Value = '63.185.76'
$ws.Range("E17").Value = '  +1.21%  '
$ws.Range("D18").Value = '2.434.41'
$ws.Range("E18").Value = '  -0.22%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '11.30'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.55%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.35'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +5.52%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '328.02'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("E22").Value = '  +1.01%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.07'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +12.34%  '
$ws.Range("E24").Value = '  +0.21%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '65.62'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.49%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '613.41'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +5.61%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.96'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +4.61%  '
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("E30").Value = '  +4.02%  '
$ws.Range("E31").Value = '  +0.24%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '8.24'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.07%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.142'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.93%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.89'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.94%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '5.19'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +7.04%  '
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("E38").Value = '  -0.62%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.42'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.82%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.77'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.07%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '147.58'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.05%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.78'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.03%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.61'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +6.79%  '
$ws.Range("E44").Value = '  -0.18%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '41.87'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.62%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '148.68'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("E47").Value = '  +2.50%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '21.20'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +3.33%  '
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("E51").Value = '  +0.18%  '
